# Update cryptos list - price and 1h-volume-change refresh (GitHub Actions run)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D ("Price") always holds text in this sheet (values like "304.60" or
# "46.510.79" are display strings, not numbers). Temporarily force text format
# while writing so COM/Excel does not silently reinterpret them as numeric and
# strip formatting (e.g. "4.00" -> 4), then restore General like the rest of the sheet.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "46.510.79"
$ws.Range("D2").NumberFormat = "General"
$ws.Range("E2").Value = "  +5.42%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.296.43"
$ws.Range("D3").NumberFormat = "General"
$ws.Range("E3").Value = "  +3.18%  "

$ws.Range("E4").Value = "  +0.03%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "304.60"
$ws.Range("D5").NumberFormat = "General"
$ws.Range("E5").Value = "  +1.87%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "101.20"
$ws.Range("D6").NumberFormat = "General"
$ws.Range("E6").Value = "  +11.53%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.568"
$ws.Range("D7").NumberFormat = "General"
$ws.Range("E7").Value = "  +1.67%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.999"
$ws.Range("D8").NumberFormat = "General"
$ws.Range("E8").Value = "  -0.01%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.524"
$ws.Range("D9").NumberFormat = "General"
$ws.Range("E9").Value = "  +6.03%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "36.52"
$ws.Range("D10").NumberFormat = "General"
$ws.Range("E10").Value = "  +9.43%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0788"
$ws.Range("D11").NumberFormat = "General"
$ws.Range("E11").Value = "  +1.11%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.42"
$ws.Range("D12").NumberFormat = "General"
$ws.Range("E12").Value = "  +6.32%  "

$ws.Range("E13").Value = "  +0.11%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.642.86"
$ws.Range("D14").NumberFormat = "General"
$ws.Range("E14").Value = "  +2.98%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.291.18"
$ws.Range("D15").NumberFormat = "General"
$ws.Range("E15").Value = "  +2.72%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "13.84"
$ws.Range("D16").NumberFormat = "General"
$ws.Range("E16").Value = "  +3.28%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.814"
$ws.Range("D17").NumberFormat = "General"
$ws.Range("E17").Value = "  +4.64%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "46.477.29"
$ws.Range("D18").NumberFormat = "General"
$ws.Range("E18").Value = "  +5.79%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.10"
$ws.Range("D19").NumberFormat = "General"
$ws.Range("E19").Value = "  +5.73%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0939"
$ws.Range("D20").NumberFormat = "General"
$ws.Range("E20").Value = "  +3.46%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.02"
$ws.Range("D21").NumberFormat = "General"
$ws.Range("E21").Value = "  +0.78%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "66.15"
$ws.Range("D22").NumberFormat = "General"
$ws.Range("E22").Value = "  +3.38%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "248.89"
$ws.Range("D23").NumberFormat = "General"

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.89"
$ws.Range("D24").NumberFormat = "General"
$ws.Range("E24").Value = "  +2.93%  "

$ws.Range("E25").Value = "  +0.07%  "

$ws.Range("E26").Value = "  +3.93%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "42.62"
$ws.Range("D27").NumberFormat = "General"
$ws.Range("E27").Value = "  +8.36%  "

$ws.Range("E28").Value = "  +1.61%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.89"
$ws.Range("D29").NumberFormat = "General"
$ws.Range("E29").Value = "  +5.37%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "20.02"
$ws.Range("D30").NumberFormat = "General"
$ws.Range("E30").Value = "  +4.19%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.82"
$ws.Range("D31").NumberFormat = "General"
$ws.Range("E31").Value = "  +12.74%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.65"
$ws.Range("D32").NumberFormat = "General"
$ws.Range("E32").Value = "  +2.69%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "147.74"
$ws.Range("D33").NumberFormat = "General"
$ws.Range("E33").Value = "  -2.56%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0795"
$ws.Range("D34").NumberFormat = "General"
$ws.Range("E34").Value = "  +3.89%  "

$ws.Range("E35").Value = "  +14.89%  "

$ws.Range("E36").Value = "  +8.53%  "

$ws.Range("E37").Value = "  +0.56%  "

$ws.Range("B38").Value = "Celestia"
$ws.Range("C38").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "16.04"
$ws.Range("D38").NumberFormat = "General"
$ws.Range("E38").Value = "  +18.61%  "

$ws.Range("B39").Value = "ARBITRUM"
$ws.Range("C39").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.77"
$ws.Range("D39").NumberFormat = "General"
$ws.Range("E39").Value = "  +5.37%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "4.00"
$ws.Range("D40").NumberFormat = "General"
$ws.Range("E40").Value = "  +10.73%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.35"
$ws.Range("D41").NumberFormat = "General"
$ws.Range("E41").Value = "  +5.49%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0302"
$ws.Range("D42").NumberFormat = "General"
$ws.Range("E42").Value = "  +0.30%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.998"
$ws.Range("D43").NumberFormat = "General"
$ws.Range("E43").Value = "  -0.08%  "

$ws.Range("E44").Value = "  +10.09%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.821.31"
$ws.Range("D45").NumberFormat = "General"
$ws.Range("E45").Value = "  +0.92%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "88.20"
$ws.Range("D46").NumberFormat = "General"
$ws.Range("E46").Value = "  +20.64%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.195"
$ws.Range("D47").NumberFormat = "General"
$ws.Range("E47").Value = "  +5.83%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "73.41"
$ws.Range("D48").NumberFormat = "General"
$ws.Range("E48").Value = "  +8.16%  "

$ws.Range("E49").Value = "  +5.68%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "96.01"
$ws.Range("D50").NumberFormat = "General"
$ws.Range("E50").Value = "  +1.57%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.519.57"
$ws.Range("D51").NumberFormat = "General"
$ws.Range("E51").Value = "  +2.94%  "
